# chore: adapt column header formatting to respective input file names (#7)
# Rename the "_old"/"_new" column-header suffixes to "_FV2404"/"_FV2410",
# freeze the header row, and wrap the used range in an Excel Table (Table1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rename header row (A1:U1): "_old" -> "_FV2404", "_new" -> "_FV2410" ---
$headers = @(
    @{ Cell = "A1"; Text = "Segmentname_FV2404" },
    @{ Cell = "B1"; Text = "Segmentgruppe_FV2404" },
    @{ Cell = "C1"; Text = "Segment_FV2404" },
    @{ Cell = "D1"; Text = "Datenelement_FV2404" },
    @{ Cell = "E1"; Text = "Segment ID_FV2404" },
    @{ Cell = "F1"; Text = "Code_FV2404" },
    @{ Cell = "G1"; Text = "Qualifier_FV2404" },
    @{ Cell = "H1"; Text = "Beschreibung_FV2404" },
    @{ Cell = "I1"; Text = "Bedingungsausdruck_FV2404" },
    @{ Cell = "J1"; Text = "Bedingung_FV2404" },
    @{ Cell = "K1"; Text = "diff" },
    @{ Cell = "L1"; Text = "Segmentname_FV2410" },
    @{ Cell = "M1"; Text = "Segmentgruppe_FV2410" },
    @{ Cell = "N1"; Text = "Segment_FV2410" },
    @{ Cell = "O1"; Text = "Datenelement_FV2410" },
    @{ Cell = "P1"; Text = "Segment ID_FV2410" },
    @{ Cell = "Q1"; Text = "Code_FV2410" },
    @{ Cell = "R1"; Text = "Qualifier_FV2410" },
    @{ Cell = "S1"; Text = "Beschreibung_FV2410" },
    @{ Cell = "T1"; Text = "Bedingungsausdruck_FV2410" },
    @{ Cell = "U1"; Text = "Bedingung_FV2410" }
)

foreach ($h in $headers) {
    $ws.Range($h.Cell).Value = $h.Text
}

# --- 2) Freeze the top (header) row ---
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- 3) Convert the used range into an Excel Table ("Table1") ---
$tableRange = $ws.Range("A1:U65")
$tbl = $ws.ListObjects.Add(1, $tableRange, $null, 1, $null)
$tbl.Name = "Table1"

Write-Output "done"
